$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02101655400301894
$ws.Range("D2").Value = 0.02481752633388368
$ws.Range("E2").Value = 0.09892296873560014
$ws.Range("F2").Value = 0.5649812616646059
$ws.Range("G2").Value = 0.002406798312569132
$ws.Range("I2").Value = 0.4206507165661506
$ws.Range("K2").Value = 1.185262995587038
$ws.Range("M2").Value = 0.3813610906196701
$ws.Range("O2").Value = 1.881629758608128
$ws.Range("C3").Value = 0.0184129870248384
$ws.Range("D3").Value = 0.02268513097104119
$ws.Range("E3").Value = 0.0951807999674088
$ws.Range("F3").Value = 0.5664004778845353
$ws.Range("G3").Value = 0.002409639797639099
$ws.Range("I3").Value = 0.4225621436446154
$ws.Range("K3").Value = 1.036240164512833
$ws.Range("M3").Value = 0.3396403185953005
$ws.Range("O3").Value = 1.900692824880821
$ws.Range("C4").Value = 0.01680738971911921
$ws.Range("D4").Value = 0.02136962688881994
$ws.Range("E4").Value = 0.09301022381620427
$ws.Range("F4").Value = 0.5678201097457247
$ws.Range("G4").Value = 0.002411476264565655
$ws.Range("I4").Value = 0.4241486593886989
$ws.Range("K4").Value = 0.9444238881191893
$ws.Range("M4").Value = 0.3140570554271775
$ws.Range("O4").Value = 1.914244101005679
$ws.Range("C5").Value = 0.01615136934086792
$ws.Range("D5").Value = 0.02083202703920506
$ws.Range("E5").Value = 0.09215733320980846
$ws.Range("F5").Value = 0.5685360810889293
$ws.Range("G5").Value = 0.002412247789600631
$ws.Range("I5").Value = 0.4248986933522687
$ws.Range("K5").Value = 0.9069308774692786
$ws.Range("M5").Value = 0.3036402373050322
$ws.Range("O5").Value = 1.92022924927079
$ws.Range("C6").Value = 0.01604233436026448
$ws.Range("D6").Value = 0.02074266826947735
$ws.Range("E6").Value = 0.09201761247393847
$ws.Range("F6").Value = 0.5686632581873923
$ws.Range("G6").Value = 0.002412377300901471
$ws.Range("I6").Value = 0.4250294789571001
$ws.Range("K6").Value = 0.9007005995603663
$ws.Range("M6").Value = 0.3019110549079542
$ws.Range("O6").Value = 1.92125099458022
$ws.Range("C7").Value = 0.01679854934557312
$ws.Range("D7").Value = 0.02136238273642732
$ws.Range("E7").Value = 0.09299859377683362
$ws.Range("F7").Value = 0.5678292095440298
$ws.Range("G7").Value = 0.002411486575803966
$ws.Range("I7").Value = 0.4241583558886504
$ws.Range("K7").Value = 0.9439185535780439
$ws.Range("M7").Value = 0.3139165356654701
$ws.Range("O7").Value = 1.914322946460928
$ws.Range("C8").Value = 0.02012031516225221
$ws.Range("D8").Value = 0.0240835867939424
$ws.Range("E8").Value = 0.09760608227603029
$ws.Range("F8").Value = 0.5653565831028047
$ws.Range("G8").Value = 0.00240775904924747
$ws.Range("I8").Value = 0.4212238942855677
$ws.Range("K8").Value = 1.133946729250681
$ws.Range("M8").Value = 0.3669689026642757
$ws.Range("O8").Value = 1.887818614125834
$ws.Range("C9").Value = 0.02657764104198179
$ws.Range("D9").Value = 0.02936918327839777
$ws.Range("E9").Value = 0.1076655713164314
$ws.Range("F9").Value = 0.5648763588663144
$ws.Range("G9").Value = 0.002401174434285493
$ws.Range("I9").Value = 0.418759948437291
$ws.Range("K9").Value = 1.504007111875467
$ws.Range("M9").Value = 0.4712686495641663
$ws.Range("O9").Value = 1.850559472265843
$ws.Range("C10").Value = 0.03128624870704755
$ws.Range("D10").Value = 0.03322006861498039
$ws.Range("E10").Value = 0.1157029126085973
$ws.Range("F10").Value = 0.5672136739184381
$ws.Range("G10").Value = 0.002396774220440914
$ws.Range("I10").Value = 0.4189763103833641
$ws.Range("K10").Value = 1.774235548894978
$ws.Range("M10").Value = 0.5480652795484247
$ws.Range("O10").Value = 1.83224723229651
$ws.Range("C11").Value = 0.03342039701824717
$ws.Range("D11").Value = 0.034964594817815
$ws.Range("E11").Value = 0.1195045597869822
$ws.Range("F11").Value = 0.5688669545433527
$ws.Range("G11").Value = 0.002394866487799201
$ws.Range("I11").Value = 0.4195192658983302
$ws.Range("K11").Value = 1.896795395917081
$ws.Range("M11").Value = 0.5830404430761149
$ws.Range("O11").Value = 1.825903855314664
$ws.Range("C12").Value = 0.0342273927904273
$ws.Range("D12").Value = 0.03562412496393819
$ws.Range("E12").Value = 0.1209654318470967
$ws.Range("F12").Value = 0.5695783019094591
$ws.Range("G12").Value = 0.002394157515843764
$ws.Range("I12").Value = 0.419789139071348
$ws.Range("K12").Value = 1.943150854014505
$ws.Range("M12").Value = 0.5962903695138948
$ws.Range("O12").Value = 1.823789073300503
$ws.Range("C13").Value = 0.03405364395666766
$ws.Range("D13").Value = 0.03548213219106344
$ws.Range("E13").Value = 0.1206498551511359
$ws.Range("F13").Value = 0.5694213004087132
$ws.Range("G13").Value = 0.002394309608757524
$ws.Range("I13").Value = 0.4197281532562158
$ws.Range("K13").Value = 1.933169872582766
$ws.Range("M13").Value = 0.593436514599361
$ws.Range("O13").Value = 1.824231725564232
$ws.Range("C14").Value = 0.03348681249136121
$ws.Range("D14").Value = 0.03501887671359327
$ws.Range("E14").Value = 0.1196243182492651
$ws.Range("F14").Value = 0.5689237657429587
$ws.Range("G14").Value = 0.002394807891067698
$ws.Range("I14").Value = 0.4195401785364794
$ws.Range("K14").Value = 1.900610208731962
$ws.Range("M14").Value = 0.5841304113386343
$ws.Range("O14").Value = 1.825724103534327
$ws.Range("C15").Value = 0.03313945976387345
$ws.Range("D15").Value = 0.03473497700223049
$ws.Range("E15").Value = 0.1189989280565911
$ws.Range("F15").Value = 0.5686301311938422
$ws.Range("G15").Value = 0.002395114852505879
$ws.Range("I15").Value = 0.4194334183527744
$ws.Range("K15").Value = 1.880659205057157
$ws.Range("M15").Value = 0.5784308778678451
$ws.Range("O15").Value = 1.826675692606898
$ws.Range("C16").Value = 0.03114661669943075
$ws.Range("D16").Value = 0.03310591018946241
$ws.Range("E16").Value = 0.1154574229626633
$ws.Range("F16").Value = 0.5671175386851957
$ws.Range("G16").Value = 0.002396900780339492
$ws.Range("I16").Value = 0.4189498015589876
$ws.Range("K16").Value = 1.766218362242626
$ws.Range("M16").Value = 0.5457803613067682
$ws.Range("O16").Value = 1.832701905968833
$ws.Range("C17").Value = 0.02992204372961282
$ws.Range("D17").Value = 0.0321046421796396
$ws.Range("E17").Value = 0.1133223254089728
$ws.Range("F17").Value = 0.5663410459643288
$ws.Range("G17").Value = 0.002398020405061221
$ws.Range("I17").Value = 0.4187672190697143
$ws.Range("K17").Value = 1.695916591164064
$ws.Range("M17").Value = 0.5257604542530459
$ws.Range("O17").Value = 1.836908835611268
$ws.Range("C18").Value = 0.02921696768969184
$ws.Range("D18").Value = 0.03152805823111038
$ws.Range("E18").Value = 0.1121079398655098
$ws.Range("F18").Value = 0.5659499375062111
$ws.Range("G18").Value = 0.002398673230320504
$ws.Range("I18").Value = 0.4187040233199042
$ws.Range("K18").Value = 1.655446330562825
$ws.Range("M18").Value = 0.5142493065813198
$ws.Range("O18").Value = 1.839515476501418
$ws.Range("C19").Value = 0.02897811605998868
$ws.Range("D19").Value = 0.03133272106550322
$ws.Range("E19").Value = 0.1116991060133969
$ws.Range("F19").Value = 0.5658270352662882
$ws.Range("G19").Value = 0.002398895787063663
$ws.Range("I19").Value = 0.4186897985807647
$ws.Range("K19").Value = 1.64173793430291
$ws.Range("M19").Value = 0.510352482803512
$ws.Range("O19").Value = 1.84043009640439
$ws.Range("C20").Value = 0.03005247784096809
$ws.Range("D20").Value = 0.03221129964811098
$ws.Range("E20").Value = 0.1135481931985538
$ws.Range("F20").Value = 0.5664179567005689
$ws.Range("G20").Value = 0.002397900304135243
$ws.Range("I20").Value = 0.4187823244828834
$ws.Range("K20").Value = 1.703403926489329
$ws.Range("M20").Value = 0.5278912181853315
$ws.Range("O20").Value = 1.83644164302703
$ws.Range("C21").Value = 0.03365333639044366
$ws.Range("D21").Value = 0.03515497581319238
$ws.Range("E21").Value = 0.119924963143589
$ws.Range("F21").Value = 0.5690675855937357
$ws.Range("G21").Value = 0.002394661169014321
$ws.Range("I21").Value = 0.4195936443204928
$ws.Range("K21").Value = 1.910175290296309
$ws.Range("M21").Value = 0.5868636889529739
$ws.Range("O21").Value = 1.82527794554332
$ws.Range("C22").Value = 0.03599992883481207
$ws.Range("D22").Value = 0.03707250142548446
$ws.Range("E22").Value = 0.1242167009134008
$ws.Range("F22").Value = 0.5712965860696002
$ws.Range("G22").Value = 0.002392622545759723
$ws.Range("I22").Value = 0.420498649194414
$ws.Range("K22").Value = 2.044988968128678
$ws.Range("M22").Value = 0.6254381267771976
$ws.Range("O22").Value = 1.819657266264215
$ws.Range("C23").Value = 0.03474813992124837
$ws.Range("D23").Value = 0.03604967481366117
$ws.Range("E23").Value = 0.1219146408997887
$ws.Range("F23").Value = 0.5700612769333802
$ws.Range("G23").Value = 0.002393703450533735
$ws.Range("I23").Value = 0.419981226891899
$ws.Range("K23").Value = 1.973066716327537
$ws.Range("M23").Value = 0.6048472990664635
$ws.Range("O23").Value = 1.822503299187161
$ws.Range("C24").Value = 0.02999351185563626
$ws.Range("D24").Value = 0.0321630827361119
$ws.Range("E24").Value = 0.1134460375621344
$ws.Range("F24").Value = 0.5663830131088829
$ws.Range("G24").Value = 0.002397954573351587
$ws.Range("I24").Value = 0.4187753652289032
$ws.Range("K24").Value = 1.700019066599396
$ws.Range("M24").Value = 0.5269279042803134
$ws.Range("O24").Value = 1.836652275231188
$ws.Range("C25").Value = 0.02483693816783727
$ws.Range("D25").Value = 0.02794487683794245
$ws.Range("E25").Value = 0.1048322674735545
$ws.Range("F25").Value = 0.5645357950868473
$ws.Range("G25").Value = 0.00240287859285901
$ws.Range("I25").Value = 0.419072100909581
$ws.Range("K25").Value = 1.404180975095073
$ws.Range("M25").Value = 0.4430238656805301
$ws.Range("O25").Value = 1.859054220872167
